$p = $ppt.ActivePresentation

# --- Slide 1: Title shape -- "Manage Job Queue\nin 2.4.1" -> "Release 2.4.1" ---
$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(5)
$tr1 = $title1.TextFrame.TextRange
# Characters 17..41 (1-based, length 25) span "Manage Job Queue" + line-break + "in 2.4.1"
$seg1 = $tr1.Characters(17, 25)
$seg1.Text = "Release 2.4.1"

# --- Slide 4: picture alt text ---
$s4 = $p.Slides.Item(4)
$pic4 = $s4.Shapes.Item(2)
$pic4.AlternativeText = "Screenshot of the Job Queue page showing list of jobs in the queue and the ability to cancel a job that has not started."

# --- Slide 5: picture alt text ---
$s5 = $p.Slides.Item(5)
$pic5 = $s5.Shapes.Item(2)
$pic5.AlternativeText = "Screenshot showing that the contents of the organism-selection drop-down box are limited to the NCBI Taxonomy."

# --- Slide 6: picture alt text ---
$s6 = $p.Slides.Item(6)
$pic6 = $s6.Shapes.Item(2)
$pic6.AlternativeText = "Screenshot of main caArray page showing the two new search criteria for experiments, which are PubMed ID and Publication Author."

# --- Slide 7: three picture alt texts ---
$s7 = $p.Slides.Item(7)
$pic7a = $s7.Shapes.Item(2)
$pic7a.AlternativeText = "Screenshot showing search by PubMed ID."
$pic7b = $s7.Shapes.Item(3)
$pic7b.AlternativeText = "Screenshot showing Search by Author."
$pic7c = $s7.Shapes.Item(4)
$pic7c.AlternativeText = "Screenshot showing results from Search by Author."

# --- Slide 8: colorize lead-in phrases of bullet paragraphs, reword last bullet ---
$s8 = $p.Slides.Item(8)
$body8 = $s8.Shapes.Item(2)
$tr8 = $body8.TextFrame.TextRange

$blue = 12611584  # RGB(0, 112, 192) encoded as BGR long used by PowerPoint COM

# Paragraph 2: "Easier upload and import of large data sets" (prefix) -> blue
$p2seg = $tr8.Characters(19, 43)
$p2seg.Font.Color.RGB = $blue

# Paragraph 3: "A plug-in architecture" (prefix) -> blue
$p3seg = $tr8.Characters(263, 22)
$p3seg.Font.Color.RGB = $blue

# Paragraph 4: "Update the" -> "An updated" (same length, 10 chars)
$p4word = $tr8.Characters(521, 10)
$p4word.Text = "An updated"

# Paragraph 4: "An updated technology " -> blue
$p4seg1 = $tr8.Characters(521, 22)
$p4seg1.Font.Color.RGB = $blue

# Paragraph 4: "stack" -> blue
$p4seg2 = $tr8.Characters(543, 5)
$p4seg2.Font.Color.RGB = $blue
